$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename Table1 column 8 header: "Total Completed" -> "codeChallenges.totalCompleted"
$ws.Range("I3").Value = "codeChallenges.totalCompleted"

# 2. Update Leaderboard rank value for Filpill (row 4)
$ws.Range("D4").Value = 23246

# 3. Update Total Attempts (H) / Total Completed (I) counters for Table2 rows 10-95
$hiData = @{
    10 = @(490952, 28871)
    11 = @(401296, 53804)
    12 = @(84579, 5878)
    13 = @(310984, 41320)
    14 = @(288336, 45636)
    15 = @(182951, 49244)
    16 = @(182323, 50129)
    17 = @(383510, 52042)
    18 = @(223173, 58590)
    19 = @(228681, 59478)
    20 = @(322294, 64683)
    21 = @(338108, 87600)
    22 = @(534993, 88980)
    23 = @(479352, 114415)
    24 = @(634541, 115582)
    25 = @(48852, 14532)
    26 = @(677513, 163337)
    27 = @(431749, 175487)
    28 = @(701670, 191546)
    29 = @(785317, 197234)
    30 = @(634698, 215163)
    31 = @(724924, 233374)
    32 = @(884476, 295661)
    33 = @(22823, 2199)
    34 = @(13235, 2246)
    35 = @(27818, 3623)
    36 = @(27201, 6029)
    37 = @(29821, 2202)
    38 = @(13036, 2665)
    39 = @(19491, 2680)
    40 = @(44099, 2878)
    41 = @(10865, 3804)
    42 = @(36518, 5810)
    43 = @(32102, 7066)
    44 = @(27258, 7104)
    45 = @(35626, 8779)
    46 = @(17878, 1765)
    47 = @(51262, 9952)
    48 = @(3574, 323)
    49 = @(42051, 4643)
    50 = @(43249, 10501)
    51 = @(56107, 10538)
    52 = @(62368, 10775)
    53 = @(37004, 15087)
    54 = @(40825, 6525)
    55 = @(30888, 2632)
    56 = @(37026, 8322)
    57 = @(19091, 2214)
    58 = @(9390, 841)
    59 = @(26977, 1151)
    60 = @(32524, 4732)
    61 = @(29129, 6053)
    62 = @(62816, 16439)
    63 = @(13294, 1296)
    64 = @(17239, 1542)
    65 = @(34363, 3276)
    66 = @(73810, 8201)
    67 = @(28298, 6627)
    68 = @(35073, 5883)
    69 = @(52025, 10471)
    70 = @(183405, 15344)
    71 = @(85460, 7492)
    72 = @(58947, 15413)
    73 = @(144254, 18319)
    74 = @(108947, 30873)
    75 = @(451018, 66011)
    76 = @(121325, 53056)
    77 = @(322492, 32343)
    78 = @(461929, 57516)
    79 = @(276711, 99165)
    80 = @(308873, 83126)
    81 = @(380283, 38661)
    82 = @(235867, 54323)
    83 = @(312191, 110410)
    84 = @(712107, 191859)
    85 = @(170158, 39823)
    86 = @(428365, 59038)
    87 = @(289120, 62044)
    88 = @(324943, 58366)
    89 = @(184118, 67621)
    90 = @(313671, 115655)
    91 = @(503209, 141069)
    92 = @(278826, 124793)
    93 = @(395439, 102364)
    94 = @(557091, 230361)
    95 = @(6447739, 5327832)
}
foreach ($row in $hiData.Keys) {
    $vals = $hiData[$row]
    $ws.Cells.Item([int]$row, 8).Value = $vals[0]
    $ws.Cells.Item([int]$row, 9).Value = $vals[1]
}

Write-Host "edit complete"